$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.738477885723114
$ws.Range("B1").Value = 1.407269716262817
$ws.Range("C1").Value = 4.557394981384277
$ws.Range("D1").Value = 1.717284083366394
$ws.Range("E1").Value = 1.185778737068176
